{"js": "// Replace each two-digit multiplication problem text with its new value.\n// Mapping preserves 1:1 correspondence between the old diff \"-\" lines and\n// the new \"+\" lines (every source string below is unique in the document).\nconst replacements = [\n  [\"22\u00d766=\", \"82\u00d774=\"],\n  [\"82\u00d798=\", \"43\u00d738=\"],\n  [\"38\u00d788=\", \"54\u00d734=\"],\n  [\"30\u00d741=\", \"86\u00d793=\"],\n  [\"73\u00d740=\", \"85\u00d750=\"],\n  [\"16\u00d787=\", \"35\u00d750=\"],\n  [\"98\u00d737=\", \"86\u00d714=\"],\n  [\"57\u00d798=\", \"95\u00d755=\"],\n  [\"59\u00d725=\", \"25\u00d767=\"],\n  [\"86\u00d782=\", \"69\u00d735=\"],\n  [\"94\u00d794=\", \"12\u00d739=\"],\n  [\"83\u00d774=\", \"28\u00d790=\"],\n  [\"11\u00d731=\", \"70\u00d749=\"],\n  [\"12\u00d757=\", \"15\u00d754=\"],\n  [\"71\u00d766=\", \"21\u00d756=\"],\n  [\"11\u00d746=\", \"18\u00d792=\"],\n  [\"68\u00d719=\", \"25\u00d796=\"],\n  [\"97\u00d739=\", \"31\u00d730=\"],\n  [\"39\u00d756=\", \"47\u00d762=\"],\n  [\"13\u00d796=\", \"84\u00d765=\"],\n  [\"42\u00d780=\", \"76\u00d782=\"],\n  [\"87\u00d766=\", \"33\u00d755=\"],\n  [\"67\u00d743=\", \"26\u00d719=\"],\n  [\"92\u00d740=\", \"41\u00d738=\"],\n  [\"74\u00d740=\", \"48\u00d720=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication problem text with its new value.\n# Mapping preserves 1:1 correspondence between the old diff \"-\" lines and\n# the new \"+\" lines (every source string below is unique in the document,\n# and none of the replacement values collide with any other source value,\n# so a simple ordered Find/ReplaceAll pass is safe).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"22\u00d766=\", \"82\u00d774=\"),\n    @(\"82\u00d798=\", \"43\u00d738=\"),\n    @(\"38\u00d788=\", \"54\u00d734=\"),\n    @(\"30\u00d741=\", \"86\u00d793=\"),\n    @(\"73\u00d740=\", \"85\u00d750=\"),\n    @(\"16\u00d787=\", \"35\u00d750=\"),\n    @(\"98\u00d737=\", \"86\u00d714=\"),\n    @(\"57\u00d798=\", \"95\u00d755=\"),\n    @(\"59\u00d725=\", \"25\u00d767=\"),\n    @(\"86\u00d782=\", \"69\u00d735=\"),\n    @(\"94\u00d794=\", \"12\u00d739=\"),\n    @(\"83\u00d774=\", \"28\u00d790=\"),\n    @(\"11\u00d731=\", \"70\u00d749=\"),\n    @(\"12\u00d757=\", \"15\u00d754=\"),\n    @(\"71\u00d766=\", \"21\u00d756=\"),\n    @(\"11\u00d746=\", \"18\u00d792=\"),\n    @(\"68\u00d719=\", \"25\u00d796=\"),\n    @(\"97\u00d739=\", \"31\u00d730=\"),\n    @(\"39\u00d756=\", \"47\u00d762=\"),\n    @(\"13\u00d796=\", \"84\u00d765=\"),\n    @(\"42\u00d780=\", \"76\u00d782=\"),\n    @(\"87\u00d766=\", \"33\u00d755=\"),\n    @(\"67\u00d743=\", \"26\u00d719=\"),\n    @(\"92\u00d740=\", \"41\u00d738=\"),\n    @(\"74\u00d740=\", \"48\u00d720=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
